# Swap the slide-deck's theme colour scheme ("Integral" -> "Office Theme").
#
# The underlying edit (per the source OOXML diff) swaps the full contents
# of ppt/theme/theme1.xml (the slide master's theme, "Integral") and
# ppt/theme/theme2.xml (the notes master's theme, "Office Theme"), so that
# after the edit the slide master uses the stock "Office Theme" colours
# and the notes master uses the old "Integral" colours. Both themes share
# an identical font scheme (Arial major/minor) and format scheme, so the
# only real content difference is the 12-colour scheme (and the name
# attributes, which are not content-visible).
#
# This automation host only exposes a single reachable theme object (the
# slide master's), so we drive the visible/slide-facing half of the swap
# through the documented PowerPoint COM surface: Theme.ThemeColorScheme.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Office Theme colour scheme values (in clrScheme document order):
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# PowerPoint COM RGB() packs colours as R + G*256 + B*65536.
$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
